$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 80
$ws.Range("H80").Value = 565
$ws.Range("I80").Value = 100
$ws.Range("J80").Value = 797.5
$ws.Range("K80").Value = 300
$ws.Range("L80").Value = 2392.5
$ws.Range("M80").Value = 698
$ws.Range("N80").Value = -4388.5

# Row 83
$ws.Range("H83").Value = 565
$ws.Range("I83").Value = 100
$ws.Range("J83").Value = 797.5
$ws.Range("K83").Value = 900
$ws.Range("L83").Value = 7177.5
$ws.Range("M83").Value = 4092
$ws.Range("N83").Value = -17161.5

# Row 88
$ws.Range("H88").Value = 2355.5557
$ws.Range("J88").Value = 2275
$ws.Range("L88").Value = 2275
$ws.Range("N88").Value = -3087

# Row 91
$ws.Range("H91").Value = 2355.5557
$ws.Range("J91").Value = 2275
$ws.Range("L91").Value = 2275
$ws.Range("N91").Value = -5083

# Row 116
$ws.Range("H116").Value = 4500
$ws.Range("I116").Value = 3000
$ws.Range("K116").Value = 3000
$ws.Range("M116").Value = 442

# Row 135
$ws.Range("H135").Value = 536.2222
$ws.Range("I135").Value = 478
$ws.Range("J135").Value = 582.8
$ws.Range("K135").Value = 4302
$ws.Range("L135").Value = 5245.2
$ws.Range("M135").Value = -1767
$ws.Range("N135").Value = -10315.2

$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 250.57143
$ws.Range("I4").Value = 307.2
$ws.Range("J4").Value = 109
$ws.Range("K4").Value = 307.2
$ws.Range("L4").Value = 109
$ws.Range("M4").Value = -191.2
$ws.Range("N4").Value = -341

$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 1794.2858
$ws.Range("I5").Value = 2002
$ws.Range("K5").Value = 2002
$ws.Range("M5").Value = -1889

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1072
$ws.Range("I16").Value = 904.63635
$ws.Range("K16").Value = 904.63635
$ws.Range("M16").Value = -617.63635

# Row 62
$ws.Range("H62").Value = 2309.4666
$ws.Range("I62").Value = 1861.9166
$ws.Range("K62").Value = 1861.9166
$ws.Range("M62").Value = -1237.9166

# Row 65
$ws.Range("H65").Value = 2309.4666
$ws.Range("I65").Value = 1861.9166
$ws.Range("K65").Value = 9309.583000000001
$ws.Range("M65").Value = -6189.583000000001

# Row 113
$ws.Range("H113").Value = 1072
$ws.Range("I113").Value = 904.63635
$ws.Range("K113").Value = 904.63635
$ws.Range("M113").Value = 1265.36365

# Row 134
$ws.Range("H134").Value = 873
$ws.Range("I134").Value = 873
$ws.Range("K134").Value = 2619
$ws.Range("M134").Value = -84

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 250.8077
$ws.Range("I2").Value = 111.933334
$ws.Range("K2").Value = 671.600004
$ws.Range("M2").Value = -558.600004

# Row 4
$ws.Range("H4").Value = 1571.619
$ws.Range("J4").Value = 2031.3334
$ws.Range("L4").Value = 6094.0002
$ws.Range("N4").Value = -6318.0002

# Row 55
$ws.Range("H55").Value = 3706.9707
$ws.Range("I55").Value = 1434
$ws.Range("K55").Value = 4302
$ws.Range("M55").Value = -4125

# Row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 45
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()

# Row 126
$ws.Range("H126").Value = 2006
$ws.Range("I126").Value = 2012
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 6036
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -3566
$ws.Range("N126").Value = -10940

# Row 132
$ws.Range("H132").Value = 1544.5
$ws.Range("I132").Value = 1471.6666
$ws.Range("J132").Value = 2200
$ws.Range("K132").Value = 4414.9998
$ws.Range("L132").Value = 6600
$ws.Range("M132").Value = -1884.9998
$ws.Range("N132").Value = -11660

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2056.6
$ws.Range("I16").Value = 2070.75
$ws.Range("K16").Value = 2070.75
$ws.Range("M16").Value = -1900.75

# Row 22
$ws.Range("H22").Value = 1164.6666
$ws.Range("I22").Value = 897
$ws.Range("J22").Value = 1700
$ws.Range("K22").Value = 897
$ws.Range("L22").Value = 1700
$ws.Range("M22").Value = -602
$ws.Range("N22").Value = -2290

# Row 27
$ws.Range("H27").Value = 1164.6666
$ws.Range("I27").Value = 897
$ws.Range("J27").Value = 1700
$ws.Range("K27").Value = 897
$ws.Range("L27").Value = 1700
$ws.Range("M27").Value = -790
$ws.Range("N27").Value = -1914

# Row 40
$ws.Range("H40").Value = 10499.5
$ws.Range("I40").Value = 999
$ws.Range("J40").Value = 20000
$ws.Range("K40").Value = 999
$ws.Range("L40").Value = 20000
$ws.Range("M40").Value = -863
$ws.Range("N40").Value = -20272

# Row 46
$ws.Range("H46").Value = 1000
$ws.Range("I46").Value = 1000
$ws.Range("K46").Value = 1000
$ws.Range("M46").Value = -812

# Row 122
$ws.Range("H122").Value = 3612.2354
$ws.Range("I122").Value = 3118
$ws.Range("K122").Value = 9354
$ws.Range("M122").Value = -6904

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 4000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 4000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -5248

# Row 65
$ws.Range("H65").Value = 4000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 20000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -26240

# Row 122
$ws.Range("H122").Value = 1106.7778
$ws.Range("I122").Value = 851.8570999999999
$ws.Range("J122").Value = 1999
$ws.Range("K122").Value = 2555.5713
$ws.Range("L122").Value = 5997
$ws.Range("M122").Value = -105.5712999999996
$ws.Range("N122").Value = -10897

# Row 126
$ws.Range("H126").Value = 3787.1177
$ws.Range("I126").Value = 2881.4546
$ws.Range("K126").Value = 8644.363799999999
$ws.Range("M126").Value = -6174.363799999999

# Row 132
$ws.Range("H132").Value = 2427.7144
$ws.Range("I132").Value = 2332.3333
$ws.Range("K132").Value = 6996.999899999999
$ws.Range("M132").Value = -4466.999899999999
